$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: snapshot current (pre-edit) values for rows 4-13, columns A,B,D,E,F,G,H,Q,R
# so the permutation can be applied without clobbering source data.
$snapshot = @{}
$snapshot["A4"] = $ws.Range("A4").Value2
$snapshot["B4"] = $ws.Range("B4").Value2
$snapshot["D4"] = $ws.Range("D4").Value2
$snapshot["E4"] = $ws.Range("E4").Value2
$snapshot["F4"] = $ws.Range("F4").Value2
$snapshot["G4"] = $ws.Range("G4").Value2
$snapshot["H4"] = $ws.Range("H4").Value2
$snapshot["Q4"] = $ws.Range("Q4").Value2
$snapshot["R4"] = $ws.Range("R4").Value2
$snapshot["A5"] = $ws.Range("A5").Value2
$snapshot["B5"] = $ws.Range("B5").Value2
$snapshot["D5"] = $ws.Range("D5").Value2
$snapshot["E5"] = $ws.Range("E5").Value2
$snapshot["F5"] = $ws.Range("F5").Value2
$snapshot["G5"] = $ws.Range("G5").Value2
$snapshot["H5"] = $ws.Range("H5").Value2
$snapshot["Q5"] = $ws.Range("Q5").Value2
$snapshot["R5"] = $ws.Range("R5").Value2
$snapshot["A6"] = $ws.Range("A6").Value2
$snapshot["B6"] = $ws.Range("B6").Value2
$snapshot["D6"] = $ws.Range("D6").Value2
$snapshot["E6"] = $ws.Range("E6").Value2
$snapshot["F6"] = $ws.Range("F6").Value2
$snapshot["G6"] = $ws.Range("G6").Value2
$snapshot["H6"] = $ws.Range("H6").Value2
$snapshot["Q6"] = $ws.Range("Q6").Value2
$snapshot["R6"] = $ws.Range("R6").Value2
$snapshot["A7"] = $ws.Range("A7").Value2
$snapshot["B7"] = $ws.Range("B7").Value2
$snapshot["D7"] = $ws.Range("D7").Value2
$snapshot["E7"] = $ws.Range("E7").Value2
$snapshot["F7"] = $ws.Range("F7").Value2
$snapshot["G7"] = $ws.Range("G7").Value2
$snapshot["H7"] = $ws.Range("H7").Value2
$snapshot["Q7"] = $ws.Range("Q7").Value2
$snapshot["R7"] = $ws.Range("R7").Value2
$snapshot["A8"] = $ws.Range("A8").Value2
$snapshot["B8"] = $ws.Range("B8").Value2
$snapshot["D8"] = $ws.Range("D8").Value2
$snapshot["E8"] = $ws.Range("E8").Value2
$snapshot["F8"] = $ws.Range("F8").Value2
$snapshot["G8"] = $ws.Range("G8").Value2
$snapshot["H8"] = $ws.Range("H8").Value2
$snapshot["Q8"] = $ws.Range("Q8").Value2
$snapshot["R8"] = $ws.Range("R8").Value2
$snapshot["A9"] = $ws.Range("A9").Value2
$snapshot["B9"] = $ws.Range("B9").Value2
$snapshot["D9"] = $ws.Range("D9").Value2
$snapshot["E9"] = $ws.Range("E9").Value2
$snapshot["F9"] = $ws.Range("F9").Value2
$snapshot["G9"] = $ws.Range("G9").Value2
$snapshot["H9"] = $ws.Range("H9").Value2
$snapshot["Q9"] = $ws.Range("Q9").Value2
$snapshot["R9"] = $ws.Range("R9").Value2
$snapshot["A10"] = $ws.Range("A10").Value2
$snapshot["B10"] = $ws.Range("B10").Value2
$snapshot["D10"] = $ws.Range("D10").Value2
$snapshot["E10"] = $ws.Range("E10").Value2
$snapshot["F10"] = $ws.Range("F10").Value2
$snapshot["G10"] = $ws.Range("G10").Value2
$snapshot["H10"] = $ws.Range("H10").Value2
$snapshot["Q10"] = $ws.Range("Q10").Value2
$snapshot["R10"] = $ws.Range("R10").Value2
$snapshot["A11"] = $ws.Range("A11").Value2
$snapshot["B11"] = $ws.Range("B11").Value2
$snapshot["D11"] = $ws.Range("D11").Value2
$snapshot["E11"] = $ws.Range("E11").Value2
$snapshot["F11"] = $ws.Range("F11").Value2
$snapshot["G11"] = $ws.Range("G11").Value2
$snapshot["H11"] = $ws.Range("H11").Value2
$snapshot["Q11"] = $ws.Range("Q11").Value2
$snapshot["R11"] = $ws.Range("R11").Value2
$snapshot["A12"] = $ws.Range("A12").Value2
$snapshot["B12"] = $ws.Range("B12").Value2
$snapshot["D12"] = $ws.Range("D12").Value2
$snapshot["E12"] = $ws.Range("E12").Value2
$snapshot["F12"] = $ws.Range("F12").Value2
$snapshot["G12"] = $ws.Range("G12").Value2
$snapshot["H12"] = $ws.Range("H12").Value2
$snapshot["Q12"] = $ws.Range("Q12").Value2
$snapshot["R12"] = $ws.Range("R12").Value2
$snapshot["A13"] = $ws.Range("A13").Value2
$snapshot["B13"] = $ws.Range("B13").Value2
$snapshot["D13"] = $ws.Range("D13").Value2
$snapshot["E13"] = $ws.Range("E13").Value2
$snapshot["F13"] = $ws.Range("F13").Value2
$snapshot["G13"] = $ws.Range("G13").Value2
$snapshot["H13"] = $ws.Range("H13").Value2
$snapshot["Q13"] = $ws.Range("Q13").Value2
$snapshot["R13"] = $ws.Range("R13").Value2

# Step 2: write the shuffled values back out per the target mapping.
$ws.Range("A4").Value = $snapshot["A7"]
$ws.Range("B4").Value = $snapshot["B7"]
$ws.Range("D4").Value = $snapshot["D7"]
$ws.Range("E4").Value = $snapshot["E7"]
$ws.Range("F4").Value = $snapshot["F7"]
$ws.Range("G4").Value = $snapshot["G7"]
$ws.Range("H4").Value = $snapshot["H7"]
$ws.Range("Q4").Value = $snapshot["Q7"]
$ws.Range("R4").Value = $snapshot["R7"]
$ws.Range("A5").Value = $snapshot["A8"]
$ws.Range("B5").Value = $snapshot["B8"]
$ws.Range("D5").Value = $snapshot["D8"]
$ws.Range("E5").Value = $snapshot["E8"]
$ws.Range("F5").Value = $snapshot["F8"]
$ws.Range("G5").Value = $snapshot["G8"]
$ws.Range("H5").Value = $snapshot["H8"]
$ws.Range("Q5").Value = $snapshot["Q8"]
$ws.Range("R5").Value = $snapshot["R8"]
$ws.Range("A6").Value = $snapshot["A9"]
$ws.Range("B6").Value = $snapshot["B9"]
$ws.Range("D6").Value = $snapshot["D9"]
$ws.Range("E6").Value = $snapshot["E9"]
$ws.Range("F6").Value = $snapshot["F9"]
$ws.Range("G6").Value = $snapshot["G9"]
$ws.Range("H6").Value = $snapshot["H9"]
$ws.Range("Q6").Value = $snapshot["Q9"]
$ws.Range("R6").Value = $snapshot["R9"]
$ws.Range("A7").Value = $snapshot["A4"]
$ws.Range("B7").Value = $snapshot["B4"]
$ws.Range("D7").Value = $snapshot["D4"]
$ws.Range("E7").Value = $snapshot["E4"]
$ws.Range("F7").Value = $snapshot["F4"]
$ws.Range("G7").Value = $snapshot["G4"]
$ws.Range("H7").Value = $snapshot["H4"]
$ws.Range("Q7").Value = $snapshot["Q4"]
$ws.Range("R7").Value = $snapshot["R4"]
$ws.Range("A8").Value = $snapshot["A5"]
$ws.Range("B8").Value = $snapshot["B5"]
$ws.Range("D8").Value = $snapshot["D5"]
$ws.Range("E8").Value = $snapshot["E5"]
$ws.Range("F8").Value = $snapshot["F5"]
$ws.Range("G8").Value = $snapshot["G5"]
$ws.Range("H8").Value = $snapshot["H5"]
$ws.Range("Q8").Value = $snapshot["Q5"]
$ws.Range("R8").Value = $snapshot["R5"]
$ws.Range("A9").Value = $snapshot["A10"]
$ws.Range("B9").Value = $snapshot["B10"]
$ws.Range("D9").Value = $snapshot["D10"]
$ws.Range("E9").Value = $snapshot["E10"]
$ws.Range("F9").Value = $snapshot["F10"]
$ws.Range("G9").Value = $snapshot["G10"]
$ws.Range("H9").Value = $snapshot["H10"]
$ws.Range("Q9").Value = $snapshot["Q10"]
$ws.Range("R9").Value = $snapshot["R10"]
$ws.Range("A10").Value = $snapshot["A11"]
$ws.Range("B10").Value = $snapshot["B11"]
$ws.Range("D10").Value = $snapshot["D11"]
$ws.Range("E10").Value = $snapshot["E11"]
$ws.Range("F10").Value = $snapshot["F11"]
$ws.Range("G10").Value = $snapshot["G11"]
$ws.Range("H10").Value = $snapshot["H11"]
$ws.Range("Q10").Value = $snapshot["Q11"]
$ws.Range("R10").Value = $snapshot["R11"]
$ws.Range("A11").Value = $snapshot["A12"]
$ws.Range("B11").Value = $snapshot["B12"]
$ws.Range("D11").Value = $snapshot["D12"]
$ws.Range("E11").Value = $snapshot["E12"]
$ws.Range("F11").Value = $snapshot["F12"]
$ws.Range("G11").Value = $snapshot["G12"]
$ws.Range("H11").Value = $snapshot["H12"]
$ws.Range("Q11").Value = $snapshot["Q12"]
$ws.Range("R11").Value = $snapshot["R12"]
$ws.Range("A12").Value = $snapshot["A13"]
$ws.Range("B12").Value = $snapshot["B13"]
$ws.Range("D12").Value = $snapshot["D13"]
$ws.Range("E12").Value = $snapshot["E13"]
$ws.Range("F12").Value = $snapshot["F13"]
$ws.Range("G12").Value = $snapshot["G13"]
$ws.Range("H12").Value = $snapshot["H13"]
$ws.Range("Q12").Value = $snapshot["Q13"]
$ws.Range("R12").Value = $snapshot["R13"]
$ws.Range("A13").Value = $snapshot["A6"]
$ws.Range("B13").Value = $snapshot["B6"]
$ws.Range("D13").Value = $snapshot["D6"]
$ws.Range("E13").Value = $snapshot["E6"]
$ws.Range("F13").Value = $snapshot["F6"]
$ws.Range("G13").Value = $snapshot["G6"]
$ws.Range("H13").Value = $snapshot["H6"]
$ws.Range("Q13").Value = $snapshot["Q6"]
$ws.Range("R13").Value = $snapshot["R6"]
